# Applies the cryptos-list price/volume refresh described in the commit.
# Numeric-looking Price (D) values are entered with a leading apostrophe so
# Excel stores them as literal text (matching the source data, which keeps
# trailing zeros / exact decimal strings like "0.150" or "2.40").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.876.84'
$ws.Range("E2").Value = '  -0.28%  '

$ws.Range("D3").Value = '3.058.50'
$ws.Range("E3").Value = '  -0.33%  '

$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = "'547.67"
$ws.Range("E5").Value = '  +1.82%  '

$ws.Range("D6").Value = "'136.71"
$ws.Range("E6").Value = '  -0.26%  '

$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = '3.056.79'
$ws.Range("E8").Value = '  -0.15%  '

$ws.Range("D9").Value = "'0.498"
$ws.Range("E9").Value = '  +1.18%  '

$ws.Range("D10").Value = "'6.30"
$ws.Range("E10").Value = '  +1.74%  '

$ws.Range("D11").Value = "'0.150"
$ws.Range("E11").Value = '  -3.73%  '

$ws.Range("D12").Value = "'0.454"
$ws.Range("E12").Value = '  +0.16%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = "'35.22"
$ws.Range("E13").Value = '  +2.38%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = "'0.0000225"
$ws.Range("E14").Value = '  +1.72%  '

$ws.Range("D15").Value = '3.537.57'
$ws.Range("E15").Value = '  -0.77%  '

$ws.Range("D16").Value = '62.810.50'
$ws.Range("E16").Value = '  -0.43%  '

$ws.Range("E17").Value = '  -2.48%  '

$ws.Range("D18").Value = '3.048.71'
$ws.Range("E18").Value = '  -0.63%  '

$ws.Range("D19").Value = "'6.74"
$ws.Range("E19").Value = '  +1.87%  '

$ws.Range("D20").Value = "'488.05"
$ws.Range("E20").Value = '  +3.96%  '

$ws.Range("D21").Value = "'13.46"
$ws.Range("E21").Value = '  +0.07%  '

$ws.Range("D22").Value = "'0.684"
$ws.Range("E22").Value = '  -1.33%  '

$ws.Range("D23").Value = "'7.13"
$ws.Range("E23").Value = '  +1.84%  '

$ws.Range("D24").Value = "'82.35"
$ws.Range("E24").Value = '  +5.10%  '

$ws.Range("D25").Value = "'12.28"
$ws.Range("E25").Value = '  +1.30%  '

$ws.Range("E26").Value = '  +0.00%  '

$ws.Range("D27").Value = "'2.73"
$ws.Range("E27").Value = '  +1.39%  '

$ws.Range("D28").Value = "'7.94"
$ws.Range("E28").Value = '  +0.87%  '

$ws.Range("E29").Value = '  -0.23%  '

$ws.Range("D30").Value = "'1.96"
$ws.Range("E30").Value = '  +4.55%  '

$ws.Range("D31").Value = "'26.13"
$ws.Range("E31").Value = '  +0.12%  '

$ws.Range("E32").Value = '  +0.09%  '

$ws.Range("D33").Value = "'5.78"
$ws.Range("E33").Value = '  +5.67%  '

$ws.Range("D34").Value = "'2.40"
$ws.Range("E34").Value = '  +4.23%  '

$ws.Range("D35").Value = "'55.68"
$ws.Range("E35").Value = '  -5.76%  '

$ws.Range("D36").Value = "'5.97"
$ws.Range("E36").Value = '  +0.56%  '

$ws.Range("D37").Value = "'459.85"
$ws.Range("E37").Value = '  -4.50%  '

$ws.Range("D38").Value = '3.184.53'
$ws.Range("E38").Value = '  -2.51%  '

$ws.Range("D39").Value = "'0.0809"
$ws.Range("E39").Value = '  +2.15%  '

$ws.Range("E40").Value = '  -0.70%  '

$ws.Range("E41").Value = '  +1.92%  '

$ws.Range("D42").Value = "'8.21"
$ws.Range("E42").Value = '  +1.13%  '

$ws.Range("E43").Value = '  -3.15%  '

$ws.Range("D44").Value = "'26.65"
$ws.Range("E44").Value = '  +5.30%  '

$ws.Range("E45").Value = '  -0.10%  '

$ws.Range("D46").Value = "'0.249"
$ws.Range("E46").Value = '  -0.42%  '

$ws.Range("E47").Value = '  +1.94%  '

$ws.Range("D48").Value = "'2.02"
$ws.Range("E48").Value = '  +0.92%  '

$ws.Range("D49").Value = "'117.12"
$ws.Range("E49").Value = '  -4.97%  '

$ws.Range("D50").Value = '0.0₃0502'
$ws.Range("E50").Value = '  -3.25%  '

$ws.Range("B51").Value = 'ThetaToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D51").Value = "'2.09"
$ws.Range("E51").Value = '  +3.89%  '
